$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($Range, $Text) {
    $Range.NumberFormat = "@"
    $Range.Value = $Text
    $Range.Style = "Normal"
}

Set-TextValue $ws.Range('D2') '57.507.88'

Set-TextValue $ws.Range('D3') '3.016.82'
Set-TextValue $ws.Range('E3') '  +0.34%  '

Set-TextValue $ws.Range('D4') '0.999'
Set-TextValue $ws.Range('E4') '  -0.13%  '

Set-TextValue $ws.Range('D5') '509.03'
Set-TextValue $ws.Range('E5') '  -0.33%  '

Set-TextValue $ws.Range('D6') '139.94'
Set-TextValue $ws.Range('E6') '  +0.35%  '

Set-TextValue $ws.Range('E7') '  +0.01%  '

Set-TextValue $ws.Range('D8') '0.434'
Set-TextValue $ws.Range('E8') '  +0.11%  '

Set-TextValue $ws.Range('D9') '7.56'
Set-TextValue $ws.Range('E9') '  -0.14%  '

Set-TextValue $ws.Range('E10') '  +1.45%  '

Set-TextValue $ws.Range('D11') '0.366'
Set-TextValue $ws.Range('E11') '  +2.87%  '

Set-TextValue $ws.Range('D12') '3.530.27'
Set-TextValue $ws.Range('E12') '  +0.27%  '

Set-TextValue $ws.Range('E13') '  +0.59%  '

Set-TextValue $ws.Range('D14') '26.40'
Set-TextValue $ws.Range('E14') '  +2.42%  '

Set-TextValue $ws.Range('E15') '  +4.24%  '

Set-TextValue $ws.Range('B16') 'Polkadot'
Set-TextValue $ws.Range('C16') 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
Set-TextValue $ws.Range('D16') '6.26'
Set-TextValue $ws.Range('E16') '  +5.03%  '

Set-TextValue $ws.Range('B17') 'WrappedBTC'
Set-TextValue $ws.Range('C17') 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
Set-TextValue $ws.Range('D17') '57.532.16'
Set-TextValue $ws.Range('E17') '  +1.28%  '

Set-TextValue $ws.Range('D18') '3.007.65'
Set-TextValue $ws.Range('E18') '  -0.03%  '

Set-TextValue $ws.Range('E19') '  +2.72%  '

Set-TextValue $ws.Range('D20') '7.95'
Set-TextValue $ws.Range('E20') '  +1.11%  '

Set-TextValue $ws.Range('D21') '329.31'
Set-TextValue $ws.Range('E21') '  -0.28%  '

Set-TextValue $ws.Range('D22') '0.999'
Set-TextValue $ws.Range('E22') '  -0.10%  '

Set-TextValue $ws.Range('D23') '5.73'
Set-TextValue $ws.Range('E23') '  -0.75%  '

Set-TextValue $ws.Range('D24') '0.500'
Set-TextValue $ws.Range('E24') '  +3.38%  '

Set-TextValue $ws.Range('D25') '64.59'
Set-TextValue $ws.Range('E25') '  +2.44%  '

Set-TextValue $ws.Range('E26') '  -4.04%  '

Set-TextValue $ws.Range('D27') '0.999'
Set-TextValue $ws.Range('E27') '  -0.16%  '

Set-TextValue $ws.Range('D28') '0.0₃0922'
Set-TextValue $ws.Range('E28') '  +0.97%  '

Set-TextValue $ws.Range('D29') '6.78'
Set-TextValue $ws.Range('E29') '  +0.79%  '

Set-TextValue $ws.Range('D30') '7.34'
Set-TextValue $ws.Range('E30') '  +3.56%  '

Set-TextValue $ws.Range('E31') '  +0.70%  '

Set-TextValue $ws.Range('E32') '  -5.50%  '

Set-TextValue $ws.Range('D33') '20.65'
Set-TextValue $ws.Range('E33') '  -0.16%  '

Set-TextValue $ws.Range('D34') '4.76'
Set-TextValue $ws.Range('E34') '  +3.68%  '

Set-TextValue $ws.Range('D35') '154.24'
Set-TextValue $ws.Range('E35') '  +0.10%  '

Set-TextValue $ws.Range('D36') '5.89'
Set-TextValue $ws.Range('E36') '  +3.80%  '

Set-TextValue $ws.Range('D37') '1.28'
Set-TextValue $ws.Range('E37') '  +0.53%  '

Set-TextValue $ws.Range('D38') '24.55'
Set-TextValue $ws.Range('E38') '  +1.36%  '

Set-TextValue $ws.Range('D39') '0.0678'
Set-TextValue $ws.Range('E39') '  -0.23%  '

Set-TextValue $ws.Range('D40') '3.048.24'
Set-TextValue $ws.Range('E40') '  +0.29%  '

Set-TextValue $ws.Range('D41') '37.65'
Set-TextValue $ws.Range('E41') '  +1.72%  '

Set-TextValue $ws.Range('D42') '3.86'
Set-TextValue $ws.Range('E42') '  +5.25%  '

Set-TextValue $ws.Range('E43') '  -0.12%  '

Set-TextValue $ws.Range('D44') '0.650'
Set-TextValue $ws.Range('E44') '  +0.04%  '

Set-TextValue $ws.Range('E45') '  -0.41%  '

Set-TextValue $ws.Range('D46') '2.224.75'
Set-TextValue $ws.Range('E46') '  -2.25%  '

Set-TextValue $ws.Range('D47') '0.988'
Set-TextValue $ws.Range('E47') '  -1.38%  '

Set-TextValue $ws.Range('D48') '6.07'
Set-TextValue $ws.Range('E48') '  +4.05%  '

Set-TextValue $ws.Range('D49') '0.0240'
Set-TextValue $ws.Range('E49') '  +0.14%  '

Set-TextValue $ws.Range('D50') '19.49'
Set-TextValue $ws.Range('E50') '  -0.02%  '

Set-TextValue $ws.Range('D51') '1.87'
Set-TextValue $ws.Range('E51') '  -5.80%  '

